$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.317.95'
$ws.Range("E2").Value = '  +4.10%  '

$ws.Range("D3").Value = '2.210.34'
$ws.Range("E3").Value = '  +2.37%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.58'
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("E6").Value = '  -0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.75'
$ws.Range("E7").Value = '  -3.84%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.399'
$ws.Range("E9").Value = '  +2.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.94'
$ws.Range("E10").Value = '  -1.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0886'
$ws.Range("E11").Value = '  +4.95%  '

$ws.Range("E12").Value = '  -0.42%  '

$ws.Range("D13").Value = '2.543.64'
$ws.Range("E13").Value = '  +2.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.50'
$ws.Range("E14").Value = '  -2.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.44'
$ws.Range("E15").Value = '  -1.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.790'
$ws.Range("E16").Value = '  -1.88%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.51'
$ws.Range("E17").Value = '  +0.43%  '

$ws.Range("D18").Value = '2.217.58'
$ws.Range("E18").Value = '  +2.97%  '

$ws.Range("D19").Value = '41.211.02'
$ws.Range("E19").Value = '  +4.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.36'
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("E21").Value = '  +5.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("E22").Value = '  +0.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.27'
$ws.Range("E23").Value = '  +8.72%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  -0.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.52'
$ws.Range("E27").Value = '  -0.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '167.29'
$ws.Range("E28").Value = '  -2.94%  '

$ws.Range("E29").Value = '  +1.12%  '

$ws.Range("E30").Value = '  -1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.82'
$ws.Range("E31").Value = '  +0.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.55'
$ws.Range("E32").Value = '  -4.86%  '

$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("E34").Value = '  +5.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.59'
$ws.Range("E35").Value = '  +0.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0620'
$ws.Range("E36").Value = '  +0.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.49'
$ws.Range("E37").Value = '  -6.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.66'
$ws.Range("E38").Value = '  -0.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.35'
$ws.Range("E39").Value = '  -2.13%  '

$ws.Range("E40").Value = '  +0.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.000235'
$ws.Range("E41").Value = '  +27.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.77'
$ws.Range("E42").Value = '  -5.72%  '

$ws.Range("E43").Value = '  +3.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.62'
$ws.Range("E44").Value = '  +9.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0976'
$ws.Range("E45").Value = '  +6.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.49'
$ws.Range("E46").Value = '  -4.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.18'
$ws.Range("E47").Value = '  -1.56%  '

$ws.Range("D48").Value = '1.461.60'
$ws.Range("E48").Value = '  -3.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.36'
$ws.Range("E49").Value = '  -7.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.78'
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.24'
$ws.Range("E51").Value = '  +4.69%  '
